$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "27.593.29"
$ws.Range("E2").Value2 = "  -1.40%  "
$ws.Range("D3").Value2 = "1.846.95"
$ws.Range("E3").Value2 = "  -1.09%  "
$ws.Range("D4").Value2 = "'1.001"
$ws.Range("E4").Value2 = "  -0.37%  "
$ws.Range("D5").Value2 = "'313.59"
$ws.Range("E5").Value2 = "  -1.54%  "
$ws.Range("E6").Value2 = "  -0.26%  "
$ws.Range("D7").Value2 = "'0.4247"
$ws.Range("E7").Value2 = "  -1.96%  "
$ws.Range("D8").Value2 = "'0.3636"
$ws.Range("E8").Value2 = "  -2.46%  "
$ws.Range("D9").Value2 = "'44.62"
$ws.Range("E9").Value2 = "  +0.40%  "
$ws.Range("D10").Value2 = "'0.07285"
$ws.Range("E10").Value2 = "  -1.86%  "
$ws.Range("D11").Value2 = "'0.8757"
$ws.Range("E11").Value2 = "  -5.88%  "
$ws.Range("D12").Value2 = "'20.72"
$ws.Range("E12").Value2 = "  -2.27%  "
$ws.Range("D13").Value2 = "1.867.89"
$ws.Range("E13").Value2 = "  -6.38%  "
$ws.Range("D14").Value2 = "'5.337"
$ws.Range("E14").Value2 = "  -1.64%  "
$ws.Range("D15").Value2 = "'6.518"
$ws.Range("E15").Value2 = "  -3.10%  "
$ws.Range("D16").Value2 = "'0.06879"
$ws.Range("E16").Value2 = "  +0.34%  "
$ws.Range("D17").Value2 = "'1.002"
$ws.Range("E17").Value2 = "  -0.26%  "
$ws.Range("D18").Value2 = "'79.16"
$ws.Range("E18").Value2 = "  -1.92%  "
$ws.Range("D19").Value2 = "'0.000008871"
$ws.Range("E19").Value2 = "  -1.63%  "
$ws.Range("E20").Value2 = "  -0.12%  "
$ws.Range("D21").Value2 = "'15.37"
$ws.Range("E21").Value2 = "  -2.34%  "
$ws.Range("D22").Value2 = "27.604.55"
$ws.Range("E22").Value2 = "  -1.34%  "
$ws.Range("E23").Value2 = "  -2.61%  "
$ws.Range("D24").Value2 = "'10.42"
$ws.Range("E24").Value2 = "  -5.29%  "
$ws.Range("D25").Value2 = "2.088.84"
$ws.Range("E25").Value2 = "  -4.17%  "
$ws.Range("D26").Value2 = "'1.983"
$ws.Range("E26").Value2 = "  -3.08%  "
$ws.Range("D27").Value2 = "'152.44"
$ws.Range("E27").Value2 = "  -0.87%  "
$ws.Range("D28").Value2 = "'18.92"
$ws.Range("E28").Value2 = "  +2.30%  "
$ws.Range("D29").Value2 = "'121.64"
$ws.Range("E29").Value2 = "  +7.71%  "
$ws.Range("D30").Value2 = "'5.256"
$ws.Range("E30").Value2 = "  -4.49%  "
$ws.Range("D31").Value2 = "'1.884"
$ws.Range("E31").Value2 = "  +11.31%  "
$ws.Range("D32").Value2 = "'0.08894"
$ws.Range("E32").Value2 = "  -0.91%  "
$ws.Range("D33").Value2 = "'0.7649"
$ws.Range("E33").Value2 = "  -5.19%  "
$ws.Range("D34").Value2 = "'4.562"
$ws.Range("E34").Value2 = "  -4.54%  "
$ws.Range("D35").Value2 = "'2.955"
$ws.Range("E35").Value2 = "  +0.02%  "
$ws.Range("D36").Value2 = "'1.099"
$ws.Range("E36").Value2 = "  -6.44%  "
$ws.Range("D37").Value2 = "'0.9997"
$ws.Range("E37").Value2 = "  -0.40%  "
$ws.Range("D38").Value2 = "'1.092"
$ws.Range("E38").Value2 = "  -2.46%  "
$ws.Range("D39").Value2 = "'0.05348"
$ws.Range("E39").Value2 = "  -2.48%  "
$ws.Range("D40").Value2 = "'0.01936"
$ws.Range("E40").Value2 = "  -1.68%  "
$ws.Range("D41").Value2 = "'2.806"
$ws.Range("E41").Value2 = "  -6.19%  "
$ws.Range("D42").Value2 = "'0.5102"
$ws.Range("E42").Value2 = "  -2.47%  "
$ws.Range("D43").Value2 = "'6.870"
$ws.Range("E43").Value2 = "  -1.87%  "
$ws.Range("D44").Value2 = "'0.1647"
$ws.Range("E44").Value2 = "  -2.03%  "
$ws.Range("D45").Value2 = "'8.279"
$ws.Range("E45").Value2 = "  -5.41%  "
$ws.Range("D46").Value2 = "'0.06530"
$ws.Range("E46").Value2 = "  -2.78%  "
$ws.Range("D47").Value2 = "'0.4749"
$ws.Range("E47").Value2 = "  -2.49%  "
$ws.Range("D48").Value2 = "'10.29"
$ws.Range("E48").Value2 = "  -2.51%  "
$ws.Range("D49").Value2 = "'104.65"
$ws.Range("E49").Value2 = "  -1.97%  "
$ws.Range("D50").Value2 = "'0.9999"
$ws.Range("E50").Value2 = "  -0.28%  "
$ws.Range("D51").Value2 = "'1.623"
$ws.Range("E51").Value2 = "  -2.87%  "
